# Fixed close and reactivate incident use-cases.
# Insert a new "Reactivate incident" use-case row above row 26 (pushing the
# remaining use cases, and the UUCW totals row, down by one), and fill in
# its data the same way the other "simple" (Weight=5) use cases are filled.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Insert a new row at 26 - shifts rows 26..37 (incl. the UUCW totals row)
# down to 27..38 and extends the B/E/F/H/I formulas and SUM ranges that
# covered the old range automatically.
$ws.Rows.Item(26).Insert()

# New use case: "Reactivate incident" - a simple (5 point) use case, same
# shape as the other simple rows (B=1 actor/simple flag, H=1 "done" flag).
$ws.Range("A26").Value = "Reactivate incident"
$ws.Range("B26").Value = 1
$ws.Range("E26").Formula = "=B26*5+C26*10+D26*15"
$ws.Range("F26").Formula = "=E26/E38"
$ws.Range("H26").Value = 1
$ws.Range("I26").Formula = "=E26*H26"

# Match the author's final selection/scroll state.
try {
    $excel.ActiveWindow.ScrollRow = 19
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("I26").Select() | Out-Null
